# Merge the three templated sections (narrative / model_section /
# benchmark_section) into a single {{ full_report }} placeholder.
$d = $word.ActiveDocument

# 1. Re-purpose the "{{ narrative }}" paragraph to hold the new placeholder.
$d.Content.Find.Execute("{{ narrative }}", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "{{ full_report }}", 2)

# 2. Remove the now-redundant "{{ model_section }}" and
#    "{{ benchmark_section }}" paragraphs entirely (text + paragraph mark).
#    Delete from the end backwards so earlier indices stay valid.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $text = $d.Paragraphs.Item($i).Range.Text
    if ($text -like "*{{ model_section }}*" -or $text -like "*{{ benchmark_section }}*") {
        $d.Paragraphs.Item($i).Range.Delete()
    }
}
